# Adapt column header formatting to respective input file names.
# "_old" -> "_FV2410", "_new" -> "_FV2504"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("A1:U1")
foreach ($cell in $headerRange.Cells) {
    $v = $cell.Value2
    if ($v -ne $null) {
        if ($v -like "*_old") {
            $cell.Value2 = ($v -replace "_old$", "_FV2410")
        } elseif ($v -like "*_new") {
            $cell.Value2 = ($v -replace "_new$", "_FV2504")
        }
    }
}

# Freeze the header row (row 1).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel table ("Table1") with an autofilter,
# matching the header row that was just renamed.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$tbl.Name = "Table1"
